$d = $word.ActiveDocument

# --- Change 1: add a new "Do we have the CT scan images?" bullet to the
# "Data Questions:" list, right after "How many patients have missing
# data and for what variables? (percentages)" and right before the
# "Research Questions:" heading paragraph. ---

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "How many patients have missing data and for what variables?*") {
        $target = $p
    }
}

$newPara = $target.Range.InsertParagraphAfter()
$idx = $target.Index + 1
$d.Paragraphs.Item($idx).Range.Text = "Do we have the CT scan images?"

# --- Change 2: add a new empty paragraph (indented, no list/heading
# formatting) right after the last paragraph in the document ("What is
# the MELD score and how is it determined? "). ---

$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRange = $d.Range($last.Range.End, $last.Range.End)
$blankXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p>'
$endRange.InsertXML($blankXml)
